$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B8").Value = "All Results experimented with 100 tours, 100 generations and 70% mutations"

$ws.Range("B8").Select()
